# Applies the "update scripts wuth new tpm" edit to the LR-pairs sheet:
# recomputed TPM-based stats for the Calca-Ramp1 pairs, plus two new
# "Sending cluster" / "Target cluster" = MuSCs rows, expanding the table
# from 5 data rows (A1:T6) to 8 data rows (A1:T9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> FAPs (Calca/Ramp1)
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Calca"
$ws.Cells.Item(2,3).Value = "Ramp1"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 1.0
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.21212
$ws.Cells.Item(2,8).Value = 0.63636
$ws.Cells.Item(2,9).Value = 0.2045603692733198
$ws.Cells.Item(2,10).Value = 0.2783689285053439
$ws.Cells.Item(2,11).Value = 2.0
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.7451976666666668
$ws.Cells.Item(2,14).Value = 2.235593
$ws.Cells.Item(2,15).Value = 0.06145846259783393
$ws.Cells.Item(2,16).Value = 0.06209960373184343
$ws.Cells.Item(2,17).Value = 0.1580713290533333
$ws.Cells.Item(2,18).Value = 1.42264196148
$ws.Cells.Item(2,19).Value = 0.01257196580398342
$ws.Cells.Item(2,20).Value = 0.01728660015143971

# Row 3: FAPs -> Inflammatory-Mac (Calca/Ramp1)
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Calca"
$ws.Cells.Item(3,3).Value = "Ramp1"
$ws.Cells.Item(3,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(3,5).Value = 1.0
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.21212
$ws.Cells.Item(3,8).Value = 0.63636
$ws.Cells.Item(3,9).Value = 0.2045603692733198
$ws.Cells.Item(3,10).Value = 0.2783689285053439
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 6.298285
$ws.Cells.Item(3,14).Value = 18.894855
$ws.Cells.Item(3,15).Value = 0.5194365608180895
$ws.Cells.Item(3,16).Value = 0.5248553775533562
$ws.Cells.Item(3,17).Value = 1.3359922142
$ws.Cells.Item(3,18).Value = 12.0239299278
$ws.Cells.Item(3,19).Value = 0.1062561346950116
$ws.Cells.Item(3,20).Value = 0.1461034290697955

# Row 4: FAPs -> MuSCs (Calca/Ramp1)
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Calca"
$ws.Cells.Item(4,3).Value = "Ramp1"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1.0
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.21212
$ws.Cells.Item(4,8).Value = 0.63636
$ws.Cells.Item(4,9).Value = 0.2045603692733198
$ws.Cells.Item(4,10).Value = 0.2783689285053439
$ws.Cells.Item(4,11).Value = 2.0
$ws.Cells.Item(4,12).Value = 1.0
$ws.Cells.Item(4,13).Value = 0.375557
$ws.Cells.Item(4,14).Value = 0.751114
$ws.Cells.Item(4,15).Value = 0.03097319928697403
$ws.Cells.Item(4,16).Value = 0.02086420996909538
$ws.Cells.Item(4,17).Value = 0.07966315084
$ws.Cells.Item(4,18).Value = 0.47797890504
$ws.Cells.Item(4,19).Value = 0.006335889083719534
$ws.Cells.Item(4,20).Value = 0.005807947773207595

# Row 5: FAPs -> Resolving-Mac (Calca/Ramp1)
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Calca"
$ws.Cells.Item(5,3).Value = "Ramp1"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 1.0
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.21212
$ws.Cells.Item(5,8).Value = 0.63636
$ws.Cells.Item(5,9).Value = 0.2045603692733198
$ws.Cells.Item(5,10).Value = 0.2783689285053439
$ws.Cells.Item(5,11).Value = 3.0
$ws.Cells.Item(5,12).Value = 1.0
$ws.Cells.Item(5,13).Value = 4.706185000000001
$ws.Cells.Item(5,14).Value = 14.118555
$ws.Cells.Item(5,15).Value = 0.3881317772971025
$ws.Cells.Item(5,16).Value = 0.3921808087457049
$ws.Cells.Item(5,17).Value = 0.9982759622000001
$ws.Cells.Item(5,18).Value = 8.9844836598
$ws.Cells.Item(5,19).Value = 0.0793963796906052
$ws.Cells.Item(5,20).Value = 0.1091709515109011

# Row 6: MuSCs -> FAPs (Calca/Ramp1)
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Calca"
$ws.Cells.Item(6,3).Value = "Ramp1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 2.0
$ws.Cells.Item(6,6).Value = 1.0
$ws.Cells.Item(6,7).Value = 0.8248355
$ws.Cells.Item(6,8).Value = 1.649671
$ws.Cells.Item(6,9).Value = 0.7954396307266801
$ws.Cells.Item(6,10).Value = 0.721631071494656
$ws.Cells.Item(6,11).Value = 2.0
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.7451976666666668
$ws.Cells.Item(6,14).Value = 2.235593
$ws.Cells.Item(6,15).Value = 0.06145846259783393
$ws.Cells.Item(6,16).Value = 0.06209960373184343
$ws.Cells.Item(6,17).Value = 0.6146654899838334
$ws.Cells.Item(6,18).Value = 3.687992939903001
$ws.Cells.Item(6,19).Value = 0.0488864967938505
$ws.Cells.Item(6,20).Value = 0.04481300358040372

# Row 7: MuSCs -> Inflammatory-Mac (Calca/Ramp1)
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Calca"
$ws.Cells.Item(7,3).Value = "Ramp1"
$ws.Cells.Item(7,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(7,5).Value = 2.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = 0.8248355
$ws.Cells.Item(7,8).Value = 1.649671
$ws.Cells.Item(7,9).Value = 0.7954396307266801
$ws.Cells.Item(7,10).Value = 0.721631071494656
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,12).Value = 1.0
$ws.Cells.Item(7,13).Value = 6.298285
$ws.Cells.Item(7,14).Value = 18.894855
$ws.Cells.Item(7,15).Value = 0.5194365608180895
$ws.Cells.Item(7,16).Value = 0.5248553775533562
$ws.Cells.Item(7,17).Value = 5.1950490571175
$ws.Cells.Item(7,18).Value = 31.170294342705
$ws.Cells.Item(7,19).Value = 0.4131804261230779
$ws.Cells.Item(7,20).Value = 0.3787519484835606

# Row 8: MuSCs -> MuSCs (Calca/Ramp1)
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Calca"
$ws.Cells.Item(8,3).Value = "Ramp1"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = 2.0
$ws.Cells.Item(8,6).Value = 1.0
$ws.Cells.Item(8,7).Value = 0.8248355
$ws.Cells.Item(8,8).Value = 1.649671
$ws.Cells.Item(8,9).Value = 0.7954396307266801
$ws.Cells.Item(8,10).Value = 0.721631071494656
$ws.Cells.Item(8,11).Value = 2.0
$ws.Cells.Item(8,12).Value = 1.0
$ws.Cells.Item(8,13).Value = 0.375557
$ws.Cells.Item(8,14).Value = 0.751114
$ws.Cells.Item(8,15).Value = 0.03097319928697403
$ws.Cells.Item(8,16).Value = 0.02086420996909538
$ws.Cells.Item(8,17).Value = 0.3097727458735
$ws.Cells.Item(8,18).Value = 1.239090983494
$ws.Cells.Item(8,19).Value = 0.0246373102032545
$ws.Cells.Item(8,20).Value = 0.01505626219588778

# Row 9: MuSCs -> Resolving-Mac (Calca/Ramp1)
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Calca"
$ws.Cells.Item(9,3).Value = "Ramp1"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 2.0
$ws.Cells.Item(9,6).Value = 1.0
$ws.Cells.Item(9,7).Value = 0.8248355
$ws.Cells.Item(9,8).Value = 1.649671
$ws.Cells.Item(9,9).Value = 0.7954396307266801
$ws.Cells.Item(9,10).Value = 0.721631071494656
$ws.Cells.Item(9,11).Value = 3.0
$ws.Cells.Item(9,12).Value = 1.0
$ws.Cells.Item(9,13).Value = 4.706185000000001
$ws.Cells.Item(9,14).Value = 14.118555
$ws.Cells.Item(9,15).Value = 0.3881317772971025
$ws.Cells.Item(9,16).Value = 0.3921808087457049
$ws.Cells.Item(9,17).Value = 3.881828457567501
$ws.Cells.Item(9,18).Value = 23.290970745405
$ws.Cells.Item(9,19).Value = 0.3087353976064972
$ws.Cells.Item(9,20).Value = 0.2830098572348038
